$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Row 18 is the "Client discount" line item; set its amount to 125.
$ws.Range("E18").Value = 125

# Update the footer "name, email address" placeholder with the RPA developer contact.
$ws.Range("A31").Value = "RPA Developer - RPADeveloper@Uipath.com"
